$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.098.42"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "3.851.74"
$ws.Range("E3").Value = "  +1.82%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "688.95"
$ws.Range("E5").Value = "  +3.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.10"
$ws.Range("E6").Value = "  +2.78%  "
$ws.Range("D7").Value = "3.851.48"
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("E11").Value = "  +5.26%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("E13").Value = "  +6.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.60"
$ws.Range("E14").Value = "  +3.04%  "
$ws.Range("D15").Value = "4.498.87"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").Value = "3.852.14"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "71.077.21"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.77"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.23"
$ws.Range("E19").Value = "  +1.15%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.08"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "489.53"
$ws.Range("E22").Value = "  +3.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.75"
$ws.Range("E24").Value = "  +2.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000147"
$ws.Range("E25").Value = "  +3.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.38"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.55"
$ws.Range("E27").Value = "  +2.34%  "
$ws.Range("E28").Value = "  +2.03%  "
$ws.Range("D29").Value = "4.004.40"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +9.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.64"
$ws.Range("E32").Value = "  +3.36%  "
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.77"
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.29"
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("D37").Value = "3.802.54"
$ws.Range("E37").Value = "  +1.73%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.39"
$ws.Range("E40").Value = "  +13.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.44"
$ws.Range("E41").Value = "  +2.27%  "
$ws.Range("E42").Value = "  +2.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.02"
$ws.Range("E43").Value = "  +5.66%  "
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "164.48"
$ws.Range("E46").Value = "  +3.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000309"
$ws.Range("E47").Value = "  +9.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.67"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.49"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.40"
$ws.Range("E51").Value = "  -1.77%  "
